$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# 1) "Capsule of joint" note: [NOTE}Encloses tendon...  ->  [NOTE]Encloses...
#    Only the run holding the stray "}" changes, the surrounding runs
#    ("[", "NOTE", "Encloses tendon of long head of biceps. ") must stay as
#    separate runs, so the text is swapped via FormattedText (which keeps
#    run boundaries) rather than a plain Find/Replace (which would coalesce
#    neighbouring runs that share the same formatting).
# --------------------------------------------------------------------------
$brace = $d.Content
$f = $brace.Find
$f.ClearFormatting()
$f.Text = "}"
$f.Forward = $true
$f.Wrap = 0
if ($f.Execute()) {
    $ft = $brace.FormattedText
    $ft.Text = "]"
    $brace.FormattedText = $ft
}

# --------------------------------------------------------------------------
# 2) Glenohumeral ligaments bullets "Superior" / "Middle" / "Inferior" move
#    from the second list level (w:ilvl=1) up to the first list level
#    (w:ilvl=0). Word's ListLevelNumber is 1-based, so ilvl=0 -> level 1.
# --------------------------------------------------------------------------
foreach ($label in @("Superior", "Middle", "Inferior")) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $para = $d.Paragraphs.Item($i)
        $text = $para.Range.Text.Trim()
        if ($text -eq $label -and $para.Range.ListFormat.ListLevelNumber -eq 2) {
            $para.Range.ListFormat.ListLevelNumber = 1
            break
        }
    }
}
